$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44299
$ws.Range("M2").Value = 80
$ws.Range("R2").Value = 'Provincia de Santiago'

# Row 3
$ws.Range("D3").Value = 44299
$ws.Range("M3").Value = 75
$ws.Range("R3").Value = 'Provincia de Santiago'

# Row 4
$ws.Range("D4").Value = 44971
$ws.Range("M4").Value = 25
$ws.Range("Q4").Value = '$/bandeja 5 kilos'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 5

# Row 5
$ws.Range("D5").Value = 44320
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 20
$ws.Range("R5").Value = 'Región Metropolitana'

# Row 6
$ws.Range("D6").Value = 44320
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 8000
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 8000
$ws.Range("Q6").Value = '$/bandeja 7 kilos'
$ws.Range("S6").Value = 1143
$ws.Range("T6").Value = 7

# Row 7
$ws.Range("D7").Value = 44322
$ws.Range("M7").Value = 45

# Row 8
$ws.Range("D8").Value = 44322
$ws.Range("M8").Value = 80

# Row 9
$ws.Range("D9").Value = 44300
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("S9").Value = 2143

# Row 10
$ws.Range("D10").Value = 44300
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("S10").Value = 1714

# Row 11
$ws.Range("D11").Value = 44301
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("S11").Value = 2000

# Row 12
$ws.Range("D12").Value = 44301
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("S12").Value = 1714

# Row 13
$ws.Range("D13").Value = 44302
$ws.Range("M13").Value = 50

# Row 14
$ws.Range("D14").Value = 44302
$ws.Range("M14").Value = 30

# Row 15
$ws.Range("D15").Value = 44292
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range("S15").Value = 2286

# Row 16
$ws.Range("D16").Value = 44292
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("S16").Value = 2143

# Row 17
$ws.Range("D17").Value = 44980
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 2286

# Row 18
$ws.Range("D18").Value = 44980
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 13000
$ws.Range("S18").Value = 1857
